$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, for the symbol-list refresh commit.
# NumberFormat is forced to "@" (Text) before the write so numeric-looking
# strings (prices, hour counters, ...) are stored as text, matching the
# original inlineStr cell type; Style is then reset to "Normal" so no
# stray cell formatting is left behind on the edited cells.
$updates = @{
    "D2" = "243.16"
    "G2" = "8"
    "D3" = "23.02"
    "G3" = "8"
    "D4" = "5.408"
    "G4" = "8"
    "D5" = "0.05971"
    "G5" = "8"
    "D6" = "3.423"
    "G6" = "8"
    "D7" = "6.503"
    "G7" = "8"
    "G8" = "8"
    "D9" = "0.9244"
    "G9" = "8"
    "D10" = "0.01108"
    "G10" = "8"
    "D11" = "0.1440"
    "G11" = "8"
    "D12" = "0.07432"
    "G12" = "8"
    "D13" = "0.03344"
    "G13" = "8"
    "D14" = "0.03086"
    "G14" = "8"
    "D15" = "0.09343"
    "G15" = "8"
    "D16" = "3.851"
    "G16" = "8"
    "D17" = "0.001589"
    "G17" = "8"
    "D18" = "0.04703"
    "G18" = "8"
    "D19" = "0.005875"
    "G19" = "8"
    "D20" = "0.001268"
    "G20" = "8"
    "D21" = "0.004849"
    "G21" = "8"
    "D22" = "0.00008005"
    "G22" = "8"
    "D23" = "3.576"
    "G23" = "8"
    "G24" = "8"
    "D25" = "0.3234"
    "G25" = "8"
    "G26" = "8"
    "D27" = "0.0002341"
    "G27" = "8"
    "G28" = "8"
    "G29" = "8"
    "G30" = "8"
    "G31" = "8"
    "G32" = "8"
    "G33" = "8"
    "G34" = "8"
    "G35" = "8"
    "G36" = "8"
    "G37" = "8"
    "G38" = "8"
    "G39" = "8"
    "D40" = "0.03947"
    "G40" = "8"
    "D41" = "0.006374"
    "G41" = "8"
    "D42" = "0.004003"
    "G42" = "8"
    "D43" = "0.1072"
    "G43" = "8"
    "D44" = "0.008900"
    "G44" = "8"
    "D45" = "0.00005188"
    "G45" = "8"
    "D46" = "0.00000000751"
    "G46" = "8"
    "D47" = "0.7005"
    "G47" = "8"
    "D48" = "0.002149"
    "G48" = "8"
    "G49" = "8"
    "D50" = "0.0002001"
    "G50" = "8"
    "G51" = "8"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
